$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: a new price record is inserted above the existing row 117,
# pushing the rows below it (117-189) down by one (to 118-190).
$ws.Rows("117").Insert()

$ws.Range("A117").Value = 11
$ws.Range("B117").Value = "Vega Monumental Concepción"
$ws.Range("C117").Value = "Bíobío"
$ws.Range("D117").Value = 44957
$ws.Range("E117").Value = 8
$ws.Range("F117").Value = 100112032
$ws.Range("G117").Value = "Zapallo italiano"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 220
$ws.Range("K117").Value = 5500
$ws.Range("L117").Value = 6000
$ws.Range("M117").Value = 5773
$ws.Range("N117").Value = "`$/caja 60 unidades"
$ws.Range("O117").Value = "Región de O'Higgins"
$ws.Range("P117").Value = 96
$ws.Range("Q117").Value = 60
$ws.Range("R117").Value = "Hortaliza"
